$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '30.759.25'
$ws.Range('E2').Value = '  +2.49%  '

$ws.Range('D3').Value = '1.692.47'
$ws.Range('E3').Value = '  +3.42%  '

$c = $ws.Range('D4')
$c.NumberFormat = '@'
$c.Value = '0.998'
$c.ClearFormats()
$ws.Range('E4').Value = '  -0.09%  '

$c = $ws.Range('D5')
$c.NumberFormat = '@'
$c.Value = '221.83'
$c.ClearFormats()
$ws.Range('E5').Value = '  +3.01%  '

$ws.Range('E6').Value = '  +0.38%  '

$c = $ws.Range('D7')
$c.NumberFormat = '@'
$c.Value = '0.998'
$c.ClearFormats()
$ws.Range('E7').Value = '  -0.10%  '

$c = $ws.Range('D8')
$c.NumberFormat = '@'
$c.Value = '31.04'
$c.ClearFormats()
$ws.Range('E8').Value = '  +3.71%  '

$ws.Range('E9').Value = '  +2.28%  '

$c = $ws.Range('D10')
$c.NumberFormat = '@'
$c.Value = '0.0628'
$c.ClearFormats()
$ws.Range('E10').Value = '  +2.34%  '

$ws.Range('E11').Value = '  -1.58%  '

$ws.Range('D12').Value = '1.936.82'
$ws.Range('E12').Value = '  +3.57%  '

$c = $ws.Range('D13')
$c.NumberFormat = '@'
$c.Value = '10.74'
$c.ClearFormats()
$ws.Range('E13').Value = '  +11.05%  '

$c = $ws.Range('D14')
$c.NumberFormat = '@'
$c.Value = '0.619'
$c.ClearFormats()
$ws.Range('E14').Value = '  +7.03%  '

$ws.Range('D15').Value = '1.689.85'
$ws.Range('E15').Value = '  +3.21%  '

$c = $ws.Range('D16')
$c.NumberFormat = '@'
$c.Value = '4.02'
$c.ClearFormats()
$ws.Range('E16').Value = '  +3.00%  '

$ws.Range('D17').Value = '30.736.83'
$ws.Range('E17').Value = '  +2.34%  '

$c = $ws.Range('D18')
$c.NumberFormat = '@'
$c.Value = '66.34'
$c.ClearFormats()
$ws.Range('E18').Value = '  +2.27%  '

$c = $ws.Range('D19')
$c.NumberFormat = '@'
$c.Value = '248.22'
$c.ClearFormats()
$ws.Range('E19').Value = '  -0.30%  '

$ws.Range('D20').Value = '0.0₃0720'
$ws.Range('E20').Value = '  +1.57%  '

$ws.Range('B22').Value = 'Avalanche'
$ws.Range('C22').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$c = $ws.Range('D22')
$c.NumberFormat = '@'
$c.Value = '10.25'
$c.ClearFormats()
$ws.Range('E22').Value = '  +5.47%  '

$ws.Range('B23').Value = 'Uniswap'
$ws.Range('C23').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$c = $ws.Range('D23')
$c.NumberFormat = '@'
$c.Value = '4.30'
$c.ClearFormats()
$ws.Range('E23').Value = '  +2.68%  '

$c = $ws.Range('D24')
$c.NumberFormat = '@'
$c.Value = '2.19'
$c.ClearFormats()
$ws.Range('E24').Value = '  +2.71%  '

$ws.Range('E25').Value = '  -1.49%  '

$c = $ws.Range('D26')
$c.NumberFormat = '@'
$c.Value = '15.97'
$c.ClearFormats()
$ws.Range('E26').Value = '  +1.42%  '

$ws.Range('E27').Value = '  +0.26%  '

$ws.Range('E28').Value = '  +1.61%  '

$c = $ws.Range('D29')
$c.NumberFormat = '@'
$c.Value = '0.998'
$c.ClearFormats()
$ws.Range('E29').Value = '  -0.07%  '

$ws.Range('E30').Value = '  +2.04%  '

$ws.Range('E31').Value = '  +1.15%  '

$c = $ws.Range('D32')
$c.NumberFormat = '@'
$c.Value = '3.49'
$c.ClearFormats()
$ws.Range('E32').Value = '  +2.97%  '

$ws.Range('B33').Value = 'Maker'
$ws.Range('C33').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D33').Value = '1.519.73'
$ws.Range('E33').Value = '  +6.01%  '

$ws.Range('B34').Value = 'InternetComputer(DFINITY)'
$ws.Range('C34').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$c = $ws.Range('D34')
$c.NumberFormat = '@'
$c.Value = '3.31'
$c.ClearFormats()
$ws.Range('E34').Value = '  +2.76%  '

$ws.Range('E35').Value = '  +4.91%  '

$ws.Range('E36').Value = '  -0.74%  '

$c = $ws.Range('D37')
$c.NumberFormat = '@'
$c.Value = '0.0180'
$c.ClearFormats()
$ws.Range('E37').Value = '  +4.70%  '

$c = $ws.Range('D38')
$c.NumberFormat = '@'
$c.Value = '79.87'
$c.ClearFormats()
$ws.Range('E38').Value = '  +7.61%  '

$ws.Range('B39').Value = 'MXToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$c = $ws.Range('D39')
$c.NumberFormat = '@'
$c.Value = '2.73'
$c.ClearFormats()
$ws.Range('E39').Value = '  -4.78%  '

$ws.Range('B40').Value = 'ImmutableX'
$ws.Range('C40').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$c = $ws.Range('D40')
$c.NumberFormat = '@'
$c.Value = '0.585'
$c.ClearFormats()
$ws.Range('E40').Value = '  +4.88%  '

$ws.Range('E41').Value = '  +1.37%  '

$c = $ws.Range('D42')
$c.NumberFormat = '@'
$c.Value = '0.857'
$c.ClearFormats()
$ws.Range('E42').Value = '  +2.26%  '

$ws.Range('E43').Value = '  +1.75%  '

$c = $ws.Range('D44')
$c.NumberFormat = '@'
$c.Value = '0.0503'
$c.ClearFormats()
$ws.Range('E44').Value = '  +0.91%  '

$ws.Range('E45').Value = '  -1.61%  '

$c = $ws.Range('D46')
$c.NumberFormat = '@'
$c.Value = '0.998'
$c.ClearFormats()
$ws.Range('E46').Value = '  -0.03%  '

$c = $ws.Range('D47')
$c.NumberFormat = '@'
$c.Value = '52.55'
$c.ClearFormats()
$ws.Range('E47').Value = '  -4.98%  '

$ws.Range('D48').Value = '1.828.17'

$c = $ws.Range('D49')
$c.NumberFormat = '@'
$c.Value = '5.44'
$c.ClearFormats()
$ws.Range('E49').Value = '  +0.33%  '

$c = $ws.Range('D50')
$c.NumberFormat = '@'
$c.Value = '95.97'
$c.ClearFormats()
$ws.Range('E50').Value = '  +6.15%  '

$ws.Range('E51').Value = '  +5.53%  '
